$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DRS data rows for the 41st game (SRH vs RCB), appended after row 166.
$newRows = @(
    @{ Row = 167; A = 41; B = "SRH"; C = "RCB"; D = 1; E = "RCB"; F = "SRH"; G = 20; H = "RCB"; I = "Nitin Menon"; J = "NM"; K = "Wide";   L = "Not Called"; M = "Not Called"; N = "Swapnil Singh"; O = "T Natarajan";    P = "Unsuccessful"; Q = "No" },
    @{ Row = 168; A = 41; B = "SRH"; C = "RCB"; D = 2; E = "SRH"; F = "RCB"; G = 5;  H = "SRH"; I = "Nitin Menon"; J = "NM"; K = "Wicket"; L = "Out";        M = "Out";        N = "AK Markram";    O = "Swapnil Singh"; P = "Unsuccessful"; Q = "No" },
    @{ Row = 169; A = 41; B = "SRH"; C = "RCB"; D = 2; E = "SRH"; F = "RCB"; G = 15; H = "RCB"; I = "Nitin Menon"; J = "NM"; K = "Wicket"; L = "Not Out";    M = "Not Out";    N = "B Kumar";       O = "KV Sharma";     P = "Unsuccessful"; Q = "No" },
    @{ Row = 170; A = 41; B = "SRH"; C = "RCB"; D = 2; E = "SRH"; F = "RCB"; G = 19; H = "RCB"; I = "Nitin Menon"; J = "NM"; K = "Wide";   L = "Called";     M = "Called";     N = "JD Unadkat";    O = "Mohammed Siraj"; P = "Unsuccessful"; Q = "No" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
}

# Match the author's final view state: scrolled so row 162 is at the top and
# O168 (the last bowler cell entered) is the active selection.
$excel.Goto($ws.Range("C162"), $true)
$ws.Range("O168").Select()
